# Regresion de DataSources de cuentas para PreProd (se modifican ambiente
# "Smoke AMBA/Interior" -> "Enero/ssurgwsoadev4" deja igual; el cambio real
# es sobre la fila 11 (cuenta de regresion) y la vista de la hoja.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mueve el scroll/top-left visible de la hoja de G1 a F1 (estado de vista).
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 6
$aw.ScrollRow = 1

# Fila 11: cuenta de regresion que pasa de "i-preproducciongestion..." a
# "preproducciongestion..." (se saca el prefijo "i-") y de Febrero a Marzo.
# Usamos un apostrofo inicial para conservar el prefijo de texto (quotePrefix)
# que ya tenia la celda A11.
$ws.Range("A11").Value = "'preproducciongestion.segurossura.com.ar"
$ws.Range("B11").Value = "https://preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("E11").Value = "RegreMarch"
$ws.Range("G11").Value = 24455773
$ws.Range("N11").Value = 303
